# Update mods data [2025-11-27 15:10:15]
# Append the new daily mod-count row (row 18) after the existing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the date-like text to be stored as literal text
# (matching the existing A2:A17 cells) instead of being auto-converted to a
# date serial number.
$ws.Range("A18").Value = "'2025/11/27"
$ws.Range("B18").Value = "逃离鸭科夫"
$ws.Range("C18").Value = 1270

# Match the centered formatting used by the other data rows (A3:C17).
$ws.Range("A18:C18").HorizontalAlignment = -4108
$ws.Range("A18:C18").VerticalAlignment = -4108
